$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: insert a new "Meta description" paragraph right after the
# Heading1 title paragraph (Paragraphs(1)).
# ------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$d.Paragraphs(2).Range.Style = "Normal"

# The source paragraph ("Play Arabian Fire for Free - Unique and Exotic
# Slot Experience", bold) used to be Paragraphs(51); after inserting one
# paragraph above it, it is now Paragraphs(52). Re-fetch it fresh (do not
# reuse a Range captured before the structural edit).
$srcBoldRange = $d.Paragraphs(52).Range
$newMetaPara = $d.Paragraphs(2)
$newMetaPara.Range.FormattedText = $srcBoldRange.FormattedText

# Rename the copied text to "Meta description" while keeping it bold.
$findRange = $d.Paragraphs(2).Range
$findRange.Find.ClearFormatting()
$null = $findRange.Find.Execute("Play Arabian Fire for Free - Unique and Exotic Slot Experience")
$findRange.Text = "Meta description"
$boldEndPos = $findRange.End

# Append the (non-bold) rest of the meta description text.
$metaRest = ": Read our review of Arabian Fire, the slot game with multiple jackpot opportunities, stunning visual design, and the chance to trigger the Loaded with Loot function. Play for free!"
$d.Paragraphs(2).Range.InsertAfter($metaRest)
$restRange = $d.Range($boldEndPos, $d.Paragraphs(2).Range.End - 1)
$restRange.Bold = 0

# ------------------------------------------------------------------
# Part 2: drop the old trailing bold "Play Arabian Fire for Free..."
# paragraph, and rewrite the trailing italic paragraph's text with the
# new image-generation prompt.
# ------------------------------------------------------------------

# After Part 1, the old bold paragraph shifted from 51 -> 52.
$d.Paragraphs(52).Range.Delete()

# The italic "Read our review..." paragraph is now Paragraphs(52).
$findRange2 = $d.Paragraphs(52).Range
$findRange2.Find.ClearFormatting()
$null = $findRange2.Find.Execute("Read our review of Arabian Fire, the slot game with multiple jackpot opportunities, stunning visual design, and the chance to trigger the Loaded with Loot function. Play for free!")
$findRange2.Text = "Create a cartoon-style feature image for Arabian Fire slot game that features a happy Maya warrior with glasses. The image should be vibrant and eye-catching, with the warrior shown holding a fire torch and standing in front of the game's reels. The background should feature sand dunes and an Arabian palace in the distance, with vibrant colors and intricate designs. The Maya warrior should be shown with a big smile, wearing a warrior outfit and holding a pair of glasses, indicating intelligence and fun-loving nature. The overall image should be vibrant, colorful, and convey the excitement and adventure of playing the Arabian Fire slot game."

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
